$d = $word.ActiveDocument

# Locate the paragraph that follows "LOQ4055..." (an empty paragraph),
# then the "Ver no Jupiter..." paragraph, then the "© 2020..." paragraph.
# These three whole paragraphs (including their paragraph marks) are removed,
# while the paragraph mark of the preceding "LOQ4055..." paragraph and the
# trailing empty paragraph before the page break are both kept.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "LOQ4055*") {
        $startPara = $d.Paragraphs.Item($i + 1)
    }
    if ($t -like "*2020*Contact: luizeleno@usp.br*") {
        $endPara = $p
    }
}

$start = $startPara.Range.Start
$end = $endPara.Range.End

$r = $d.Range($start, $end)
$r.Delete()
